$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.878.78"
$ws.Range("D2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.620.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "
# Row 4
$ws.Range("E4").Value = "  -0.12%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "
# Row 7
$ws.Range("E7").Value = "  +0.02%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.619.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.65%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.25%  "
# Row 11
$ws.Range("E11").Value = "  -0.66%  "
# Row 12
$ws.Range("E12").Value = "  -0.32%  "
# Row 13
$ws.Range("E13").Value = "  -1.91%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "
# Row 15
$ws.Range("E15").Value = "  +2.62%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.094.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.68%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.814.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.628.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "371.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.08%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.50%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.14%  "
# Row 24
$ws.Range("E24").Value = "  -4.25%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.15%  "
# Row 26
$ws.Range("E26").Value = "  -0.04%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.52%  "
# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.79%  "
# Row 29
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.748.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.58%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "575.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.84%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.49%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.06%  "
# Row 34
$ws.Range("E34").Value = "  -0.71%  "
# Row 35
$ws.Range("E35").Value = "  -0.01%  "
# Row 36
$ws.Range("E36").Value = "  -2.17%  "
# Row 37
$ws.Range("E37").Value = "  -1.94%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.95"
$ws.Range("D38").Style = "Normal"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.49%  "
# Row 41
$ws.Range("E41").Value = "  -0.71%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.74%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.09"
$ws.Range("D44").Style = "Normal"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0306"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.10%  "
# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "155.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.50%  "
# Row 50
$ws.Range("E50").Value = "  -2.62%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0778"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
